# "Journal de Bord" grid: headers, two log entries, and a block of empty
# styled rows ready for future entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old placeholder cell (the lone "$" in D1) -----------------
$ws.Range("D1").ClearContents()

# --- Values ---------------------------------------------------------------
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Evénement"

$ws.Range("A2").Value = "03/13/2020"
$ws.Range("B2").Value = "Fermeture des écoles et début du confinement"

$ws.Range("A3").Value = "03/18/2020"
$ws.Range("B3").Value = "Découverte des gestions de fichiers sur C, comment les écrire, lire ainsi qu'append."

# --- Column width / row heights -------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 64.05

$ws.Rows.Item(1).RowHeight = 18
$ws.Rows.Item(2).RowHeight = 15.6
$ws.Rows.Item(3).RowHeight = 31.2

# --- Header row formatting (A1:B1) -----------------------------------------
$header = $ws.Range("A1:B1")
$header.Font.Size = 14
$header.Font.Underline = $true
$header.Interior.ThemeColor = 5

$ws.Range("A1").Borders.Item(7).LineStyle = 1
$ws.Range("A1").Borders.Item(8).LineStyle = 1
$ws.Range("A1").Borders.Item(9).LineStyle = 1

$ws.Range("B1").Borders.Item(10).LineStyle = 1
$ws.Range("B1").Borders.Item(8).LineStyle = 1
$ws.Range("B1").Borders.Item(9).LineStyle = 1

$ws.Range("B1").WrapText = $true

# --- Data rows formatting (A2:B3) ------------------------------------------
$data = $ws.Range("A2:B3")
$data.Font.Size = 12
$data.Borders.LineStyle = 1

$ws.Range("B2:B3").WrapText = $true

# --- Empty styled rows below, ready for future log entries -----------------
$ws.Range("B4:B10").WrapText = $true

# --- Selection --------------------------------------------------------------
$ws.Range("B6").Select()
